# "Cleaned NA Values + Reset Working Environment"
#
# The sheet "Column Selection" lists, for each Position (FW/MF/DF/GK), the
# set of stat Columns available for that position. Two cleanup edits are
# made to the data:
#
#  1. Remove the erroneous "MF | Standard Dist" row (NA / not applicable
#     for midfielders) - this shifts every row below it up by one.
#  2. Fix the mislabeled last row for GK: "Penalty Kicks Save%" is renamed
#     to "Penalty Kicks PKA" (the correct column name).
#
# Finally, the view/selection is reset to the bottom of the refreshed list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the "MF | Standard Dist" row (row 35) - a stray/NA data row.
$ws.Rows(35).Delete()

# 2) Correct the final GK row label from "Penalty Kicks Save%" to
#    "Penalty Kicks PKA" (now sits at B100 after the row shift above).
$ws.Range("B100").Value = "Penalty Kicks PKA"

# Reset the working view: scroll near the bottom of the list and select
# just past the last populated cell, mirroring a fresh look at the sheet.
try {
    $excel.ActiveWindow.ScrollRow = 86
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("B101").Select()
